$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.189922480620155
$ws.Range("C2").Value = 0.5658914728682171
$ws.Range("J2").Value = 0.01937984496124031
$ws.Range("P2").Value = 0.1434108527131783
$ws.Range("S2").Value = 0.08139534883720931

$ws.Range("B3").Value = 0.006711409395973154
$ws.Range("C3").Value = 0.04026845637583892
$ws.Range("J3").Value = 0.06040268456375839
$ws.Range("P3").Value = 0.7181208053691275
$ws.Range("S3").Value = 0.174496644295302

$ws.Range("J4").Value = 0.04761904761904762
$ws.Range("P4").Value = 0.7380952380952381
$ws.Range("S4").Value = 0.2142857142857143

$ws.Range("B6").Value = 0.0371900826446281
$ws.Range("D6").Value = 0.008264462809917356
$ws.Range("F6").Value = 0.04545454545454546
$ws.Range("J6").Value = 0.2644628099173554
$ws.Range("O6").Value = 0.02066115702479339
$ws.Range("Q6").Value = 0.1611570247933884
$ws.Range("R6").Value = 0.0743801652892562
$ws.Range("S6").Value = 0.3884297520661157

$ws.Range("B7").Value = 0.1037735849056604
$ws.Range("D7").Value = 0.01886792452830189
$ws.Range("F7").Value = 0.04716981132075472
$ws.Range("J7").Value = 0.1226415094339623
$ws.Range("O7").Value = 0.02358490566037736
$ws.Range("Q7").Value = 0.1886792452830189
$ws.Range("R7").Value = 0.07075471698113207
$ws.Range("S7").Value = 0.4245283018867925

$ws.Range("B8").Value = 0.07335907335907337
$ws.Range("D8").Value = 0.01737451737451737
$ws.Range("F8").Value = 0.08494208494208494
$ws.Range("J8").Value = 0.1003861003861004
$ws.Range("O8").Value = 0.01158301158301158
$ws.Range("Q8").Value = 0.1853281853281853
$ws.Range("R8").Value = 0.1254826254826255
$ws.Range("S8").Value = 0.4015444015444015

$ws.Range("B9").Value = 0.08426966292134831
$ws.Range("D9").Value = 0.01123595505617977
$ws.Range("E9").Value = 0.005617977528089887
$ws.Range("F9").Value = 0.05617977528089887
$ws.Range("J9").Value = 0.1123595505617977
$ws.Range("O9").Value = 0.01123595505617977
$ws.Range("Q9").Value = 0.1853932584269663
$ws.Range("R9").Value = 0.1179775280898876
$ws.Range("S9").Value = 0.4157303370786517

$ws.Range("B10").Value = 0.09575289575289575
$ws.Range("D10").Value = 0.02007722007722008
$ws.Range("E10").Value = 0.001544401544401544
$ws.Range("F10").Value = 0.06872586872586872
$ws.Range("J10").Value = 0.1181467181467181
$ws.Range("O10").Value = 0.02162162162162162
$ws.Range("Q10").Value = 0.2030888030888031
$ws.Range("R10").Value = 0.09575289575289575
$ws.Range("S10").Value = 0.3752895752895753

$ws.Range("G11").Value = 0.1630094043887147
$ws.Range("J11").Value = 0.05015673981191222
$ws.Range("K11").Value = 0.2100313479623825
$ws.Range("L11").Value = 0.5611285266457681
$ws.Range("S11").Value = 0.01567398119122257

$ws.Range("G12").Value = 0.7292817679558011
$ws.Range("J12").Value = 0.2099447513812155
$ws.Range("L12").Value = 0.01104972375690608
$ws.Range("S12").Value = 0.04972375690607735

$ws.Range("G13").Value = 0.6956521739130435
$ws.Range("J13").Value = 0.2391304347826087
$ws.Range("S13").Value = 0.06521739130434782

$ws.Range("F15").Value = 0.0182648401826484
$ws.Range("H15").Value = 0.1917808219178082
$ws.Range("I15").Value = 0.0410958904109589
$ws.Range("J15").Value = 0.3972602739726027
$ws.Range("K15").Value = 0.0547945205479452
$ws.Range("M15").Value = 0.0091324200913242
$ws.Range("O15").Value = 0.045662100456621
$ws.Range("S15").Value = 0.2420091324200913

$ws.Range("F16").Value = 0.03571428571428571
$ws.Range("H16").Value = 0.1785714285714286
$ws.Range("I16").Value = 0.09523809523809523
$ws.Range("J16").Value = 0.3511904761904762
$ws.Range("K16").Value = 0.1666666666666667
$ws.Range("M16").Value = 0.02380952380952381
$ws.Range("O16").Value = 0.03571428571428571
$ws.Range("S16").Value = 0.1130952380952381

$ws.Range("F17").Value = 0.01279317697228145
$ws.Range("H17").Value = 0.2025586353944563
$ws.Range("I17").Value = 0.09168443496801706
$ws.Range("J17").Value = 0.4200426439232409
$ws.Range("K17").Value = 0.1044776119402985
$ws.Range("M17").Value = 0.01279317697228145
$ws.Range("O17").Value = 0.06183368869936034
$ws.Range("S17").Value = 0.09381663113006397

$ws.Range("F18").Value = 0.02459016393442623
$ws.Range("H18").Value = 0.2049180327868853
$ws.Range("I18").Value = 0.0778688524590164
$ws.Range("J18").Value = 0.4426229508196721
$ws.Range("K18").Value = 0.08196721311475409
$ws.Range("O18").Value = 0.05737704918032787
$ws.Range("S18").Value = 0.110655737704918

$ws.Range("F19").Value = 0.02118003025718608
$ws.Range("H19").Value = 0.2284417549167927
$ws.Range("I19").Value = 0.06959152798789713
$ws.Range("J19").Value = 0.3517397881996974
$ws.Range("K19").Value = 0.1043872919818457
$ws.Range("M19").Value = 0.02723146747352496
$ws.Range("N19").Value = 0.0007564296520423601
$ws.Range("O19").Value = 0.06732223903177005
$ws.Range("S19").Value = 0.1293494704992436
